# Update "想去人数" (F3) and "最低票价" (G4) values on the "展览" and
# "全部类型" sheets to reflect the newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 17
    $ws.Range("G4").Value = 45
}
